# LeaveBalanceTests.xlsx - "carry forward failure fixes"
#
# The LeaveBalance_42 .. LeaveBalance_82 scenarios (sheet rows 43-83) were
# previously excluded from the run (RunMode = "no", rows hidden by the
# RunMode autofilter). Carry them forward into the active run: flip
# RunMode to "Yes" and unhide the rows, then drop the autofilter criteria
# that was hiding them so the sheet shows everything again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LeaveBalance")

# Re-enable these 41 previously-skipped test rows. Unhide before writing
# the cell value - doing it the other way round makes the engine recompute
# (and corrupt) the wrapped-text row height.
for ($r = 43; $r -le 83; $r++) {
    $ws.Rows.Item($r).Hidden = $false
    $ws.Range("C" + $r).Value = "Yes"
}

# Clear the "RunMode = No" criterion on column 10 (RunMode) so the
# autofilter dropdown no longer hides any rows, while keeping the
# autofilter itself active over A1:V124.
$ws.Range("A1:V124").AutoFilter(10)

# Move the selection down to where the freshly unhidden rows are.
$ws.Range("C42:C124").Select()
